# NpRv VoucherCloud.docx - switch the Barcode_NpRvVoucher SDT from the
# IDAUtomationC39M (Code 39) barcode font to IDAUtomationC128M (Code 128),
# dropping the bold/caps/spacing/position/sz overrides that were only
# needed for the Code 39 rendering and tightening szCs to 20.
#
# The barcode sdt's live inside a floating text box (VML/DrawingML text
# frame), so they are not reachable through $d.Content / Find the normal
# way - we go through the Shape's TextFrame.TextRange, which exposes the
# real OOXML via WordOpenXML / InsertXML.

$d = $word.ActiveDocument

$oldRPr = '<w:rFonts w:ascii="IDAUtomationC39M" w:hAnsi="IDAUtomationC39M" w:cs="Segoe UI Light"/><w:b/><w:bCs/><w:caps/><w:color w:val="000000" w:themeColor="text1"/><w:spacing w:val="20"/><w:position w:val="20"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="sr-Latn-RS"/>'
$newRPr = '<w:rFonts w:ascii="IDAUtomationC128M" w:hAnsi="IDAUtomationC128M" w:cs="Segoe UI Light"/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:szCs w:val="20"/><w:lang w:val="sr-Latn-RS"/>'

for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    $shape = $d.Shapes.Item($i)
    $textRange = $shape.TextFrame.TextRange
    $xml = $textRange.WordOpenXML
    # NB: every shape's WordOpenXML carries the whole document font table,
    # so a bare "IDAUtomationC39M" substring check would false-positive on
    # every shape. Match on the exact old run-properties block instead,
    # which only exists in the shape that actually hosts the barcode sdt's.
    if ($xml.IndexOf($oldRPr) -ge 0) {
        $newXml = $xml.Replace($oldRPr, $newRPr)
        [void]$textRange.InsertXML($newXml)
    }
}
